$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new mission row: "Offshore Onslaught" with an "int" (intermediate) difficulty flag
$ws.Range("B6").Value = "Offshore Onslaught"
$ws.Range("C6").Value = 1

# Reflect the final selection location left by the editor
$ws.Range("E9").Select()
